$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old "blank separator" row 57,
# shifting it (and the three summary rows below it) down by one.
# The newly inserted row 57 inherits formatting from row 56 (the
# previous blank separator), leaving D57/E57/F57 blank with the
# correct number formats - exactly the row 56 pattern from before.
$ws.Rows("57").Insert()

# Fill in the new data row 56 with the extra working-hours entry.
$ws.Range("A56").Value = 2014
$ws.Range("B56").Value = 3
$ws.Range("C56").Value = 10
$ws.Range("D56").Value = 0.34375
$ws.Range("E56").Value = 0.5
$ws.Range("F56").Formula = "=(E56-D56)*24*60"
$ws.Range("G56").Formula = "=F56/60"

# Fix up the summary formulas that now live one row lower (58/59/60)
# so their ranges/refs point at the right cells after the insert.
$ws.Range("F58").Formula = "=SUM(F2:F57)"
$ws.Range("F59").Formula = "=F58/60"
$ws.Range("F60").Formula = "=F59/38.5"

# Match the workbook's recorded selection after the edit.
$ws.Range("F56").Select() | Out-Null
